$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column G (header "K") with the new Strike/K values for rows 2-11
$values = @{
    2  = 3
    3  = 10
    4  = 8
    5  = 7
    6  = 8
    7  = 4
    8  = 6
    9  = 7
    10 = 2
    11 = 2
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 7).Value = $values[$row]
}
